$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Replace the .47uF 0603 capacitor (CL10F474ZB8NNNC) with a TDK part
# --- that has a 0.55mm profile (CGB3B3X5R1E474K055AB). The old part's row
# --- is kept for reference but struck through; the new part is added as a
# --- new row directly below it.

# Mark the old (row 7) capacitor line as struck-through / discontinued
$ws.Range("A7:I7").Font.Strikethrough = $true
$ws.Range("I7").ClearContents()

# Add the replacement part on row 8, re-using row 7's schematic ref columns
$ws.Range("A8").Value = $ws.Range("A7").Value()
$ws.Range("B8").Value = $ws.Range("B7").Value()
$ws.Range("C8").Value = "CGB3B3X5R1E474K055AB"
$ws.Range("D8").Value = "445-13240-1-ND"
$ws.Range("E8").Value = "https://product.tdk.com/info/en/catalog/spec/mlccspec_commercial_lowprofile_en.pdf"
$ws.Range("F8").Value = "https://www.digikey.ca/product-detail/en/tdk-corporation/CGB3B3X5R1E474K055AB/445-13240-1-ND/3954906"
$ws.Range("G8").Value = 0.211
$ws.Range("H8").Value = 4
$ws.Range("I8").Formula = "=SUMPRODUCT(G8,H8)"

# Copy the formatting used by the other data rows onto the new row 8 cells
$ws.Range("C4").Copy()
$ws.Range("C8").PasteSpecial(-4122)
$ws.Range("D4").Copy()
$ws.Range("D8").PasteSpecial(-4122)
$ws.Range("F3").Copy()
$ws.Range("E8").PasteSpecial(-4122)
$ws.Range("F8").PasteSpecial(-4122)
$ws.Range("G4").Copy()
$ws.Range("G8").PasteSpecial(-4122)
$ws.Range("I4").Copy()
$ws.Range("I8").PasteSpecial(-4122)

# Datasheet + Digikey links for the new part
$ws.Hyperlinks.Add($ws.Range("E8"), "https://product.tdk.com/info/en/catalog/spec/mlccspec_commercial_lowprofile_en.pdf")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.digikey.ca/product-detail/en/tdk-corporation/CGB3B3X5R1E474K055AB/445-13240-1-ND/3954906")

$ws.Range("I14").Select()
